# Updated timesheet nov16 and 17
# Appends three new timesheet entries (rows 13-15) below the existing
# data on Sheet1, matching the formatting of the preceding rows.
#
# Note: this runtime's PowerShell parameter binder does not reliably
# bind *named* (-Foo bar) arguments to function parameters, so the
# helper below takes plain positional parameters instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122 ; used so the new rows inherit the same
# style (font/fill/border/number-format) as the last existing data row.
$xlPasteFormats = -4122

function Add-TimesheetRow($RowNumber, $Name, $Start, $End) {
    $srcRow = "A" + ($RowNumber - 1) + ":E" + ($RowNumber - 1)
    $dstRow = "A" + $RowNumber + ":E" + $RowNumber

    $ws.Range($srcRow).Copy()
    $ws.Range($dstRow).PasteSpecial($xlPasteFormats)
    $ws.Rows($RowNumber).RowHeight = 13.65

    $ws.Range("A" + $RowNumber).Value = $Name
    $ws.Range("B" + $RowNumber).Value = $Start
    $ws.Range("C" + $RowNumber).Value = $End
}

Add-TimesheetRow 13 "Renee Sajedian"  42690.791666666664 42690.875
Add-TimesheetRow 14 "Sophia Wilhelmi" 42690.791666666664 42690.875
Add-TimesheetRow 15 "Sophia Wilhelmi" 42691.458333333336 42691.483333333330
